# The target diff only re-orders XML attributes (e.g. <w:pgSz w:w=".." w:h=".."/>
# becomes <w:pgSz w:h=".." w:w=".."/>, xmlns:* declarations on <w:document> get
# alphabetized, <w:lsdException>/<w:style>/<w:rFonts>/<w:lang>/... attributes get
# alphabetized, etc.) across word/document.xml and word/styles.xml. Every
# attribute name/value pair is preserved verbatim - only the serialization
# order changes. That matches the commit message ("Fixed POI packaging and
# upgraded to POI 3.15"): the expected-generation fixture was simply
# regenerated by a newer Apache POI (its XMLBeans-based writer emits
# attributes in alphabetical order), with no edit to the document's text,
# formatting, styles or structure.
#
# There is no text, paragraph, run, style, section, or page-setup property
# that differs between the before/after state - so there is nothing to
# change through the Word object model. (Word's COM API - real or emulated -
# doesn't expose raw control over attribute serialization order; that's an
# internal detail of whichever XML writer produced the package.) Touching
# PageSetup/Styles values here would only risk introducing spurious diffs,
# so this script intentionally leaves the document content untouched.

$d = $word.ActiveDocument
